$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "audioFalse" condition column header becomes "currentPhase"
$ws.Range("C1").Value = "currentPhase"

# Both data rows under that column now just record the current training phase
$ws.Range("C2").Value = "train2P2"
$ws.Range("C3").Value = "train2P2"
